# Updated legacy GSC export data.
# The exported date window rolled forward by one day: the oldest day
# (2025-10-08) drops off and every remaining day's coverage numbers shift
# up one row, with no new trailing row being appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Chart")

# Remove the obsolete first data row (2025-10-08) and shift every
# subsequent row up by one, which is exactly what happened across the
# whole table (dates, not indexed / indexed / impressions all moved up
# one row, and the trailing row disappeared).
$ws.Rows("2:2").Delete()
